# Exceptions.xlsx: append three new exception rows to "Business Exceptions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$magicText = "The Conclusion Evidence Location path: \\MainFolder\Remediation_or_Justification Evidence\1-22-2020\CHR0000391114.pdf for the change: CHANGES - SOX Audit Report for magic_qq_23455.txt_07.01.73.eml made on 1/22/2020 is not a valid path."
$newText   = "The Conclusion Evidence Location path: \\MainFolder\Remediation_or_Justification Evidence\1-09-2020\CHR0000928476.pdf for the change: CHANGES - SOX Audit Report for testps9023.txt_07.01.73.eml made on 1/9/2020 is not a valid path."

# New rows 4-6 carry the same values/format used in rows 1-3.
$ws.Range("A4").Value = $magicText
$ws.Range("A5").Value = $newText
$ws.Range("A6").Value = $newText

# Match the formatting (style) that A1:A3 already use.
$ws.Range("A1").Copy()
$ws.Range("A4:A6").PasteSpecial(-4122)
$excel.CutCopyMode = $false
